# Generate Report for Archive
# Update status from "Ready for handoff" to "In Translation" for files that
# have since moved on in the localization pipeline, across all sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the status text ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E5").Value = "In Translation"
$wsOverview.Range("F5").Value = "In Translation"
$wsOverview.Range("E6").Value = "In Translation"
$wsOverview.Range("F6").Value = "In Translation"
$wsOverview.Range("E7").Value = "In Translation"
$wsOverview.Range("F7").Value = "In Translation"

# --- zh-cn sheet: column C holds the Status value ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C5").Value = "In Translation"
$wsZhCn.Range("C6").Value = "In Translation"
$wsZhCn.Range("C7").Value = "In Translation"

# --- de-de sheet: column C holds the Status value ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C5").Value = "In Translation"
$wsDeDe.Range("C6").Value = "In Translation"
$wsDeDe.Range("C7").Value = "In Translation"

# Autofit the affected columns so the column widths shrink to match the
# shorter "In Translation" text, mirroring Excel's automatic behavior.
$wsOverview.Range("E:F").Columns.AutoFit() | Out-Null
$wsZhCn.Range("C:C").Columns.AutoFit() | Out-Null
$wsDeDe.Range("C:C").Columns.AutoFit() | Out-Null
